# Updated cryptos list on Thu Nov 30 11:28:04 UTC 2023 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) columns with the latest scraped
# values, and fixes the ranking order of two coin pairs (Hedera /
# InternetComputer(DFINITY) at rows 34-35, and TrustWalletToken / FTXToken
# at rows 46-47) whose relative order changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.739.53'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.034.44'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.18'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.608'
$ws.Range('E6').Value = '  -1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.24'
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.377'
$ws.Range('E9').Value = '  -2.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0819'
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.64'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.334.98'
$ws.Range('E13').Value = '  -0.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.02'
$ws.Range('E14').Value = '  -0.28%  '
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('E16').Value = '  -2.52%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.043.50'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.688.83'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('E19').Value = '  -0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.87'
$ws.Range('E20').Value = '  -7.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.47'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('E25').Value = '  +2.76%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.38'
$ws.Range('E26').Value = '  +1.79%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.79'
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('E30').Value = '  -2.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.121'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('E32').Value = '  +8.94%  '
$ws.Range('E33').Value = '  -3.26%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.50'
$ws.Range('E34').Value = '  -1.91%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0604'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.44'
$ws.Range('E36').Value = '  +0.97%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.32'
$ws.Range('E37').Value = '  +1.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.41'
$ws.Range('E38').Value = '  +3.95%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.95'
$ws.Range('E40').Value = '  +5.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.536.25'
$ws.Range('E41').Value = '  +1.21%  '
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '96.11'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0911'
$ws.Range('E45').Value = '  -1.21%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.08'
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.11'
$ws.Range('E47').Value = '  -1.74%  '
$ws.Range('E48').Value = '  -0.68%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.11'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.224.37'
$ws.Range('E51').Value = '  -0.99%  '
